# Script: re-apply the scraper's freshly-fetched odds/rows to the
# "copa-de-la-liga-profesional" 2023 sheet.
#
# 1) A handful of existing match rows had their home/away order (and all
#    the odds/date columns that go with it, F:V) swapped between two
#    adjacent rows sharing the same kickoff date (A:E stay untouched).
# 2) Four brand-new matches (rows 131-134) are appended at the bottom.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V")

function Swap-Rows($rowA, $rowB) {
    foreach ($col in $cols) {
        $addrA = "$col$rowA"
        $addrB = "$col$rowB"
        $valA = $ws.Range($addrA).Value2
        $valB = $ws.Range($addrB).Value2
        $ws.Range($addrA).Value2 = $valB
        $ws.Range($addrB).Value2 = $valA
    }
}

# Row pairs whose F:V content got swapped.
Swap-Rows 19 20
Swap-Rows 37 38
Swap-Rows 62 63
Swap-Rows 64 65
Swap-Rows 66 67
Swap-Rows 96 97
Swap-Rows 98 99
Swap-Rows 102 103

# New rows appended at the bottom (131-134), cloning the formatting
# (border/bold/center for column A, date-time number format for column E)
# from the last existing data row so no new cell styles are introduced.
# Column D ("temporada") is the text "2023" (not the number), so the new
# cells are pre-formatted as Text before the value is written - otherwise
# a pure-digit string would silently become a number on write, same as in
# real Excel.
$ws.Range("D131:D134").NumberFormat = "@"

function Add-MatchRow {
    param(
        $row, $idx, $date,
        $home, $homeGoals, $away, $awayGoals,
        $homeOpenOdds, $homeOpenDt, $homeCloseOdds, $homeCloseDt,
        $drawOpenOdds, $drawOpenDt, $drawCloseOdds, $drawCloseDt,
        $awayOpenOdds, $awayOpenDt, $awayCloseOdds, $awayCloseDt,
        $url
    )

    $ws.Range("A$row").Value2 = $idx
    $ws.Range("A2").Copy()
    $ws.Range("A$row").PasteSpecial(-4122)

    $ws.Range("B$row").Value2 = "argentina"
    $ws.Range("C$row").Value2 = "copa-de-la-liga-profesional"
    $ws.Range("D$row").Value2 = "2023"

    $ws.Range("E$row").Value2 = $date
    $ws.Range("E2").Copy()
    $ws.Range("E$row").PasteSpecial(-4122)

    $ws.Range("F$row").Value2 = $home
    $ws.Range("G$row").Value2 = $homeGoals
    $ws.Range("H$row").Value2 = $away
    $ws.Range("I$row").Value2 = $awayGoals

    $ws.Range("J$row").Value2 = $homeOpenOdds
    $ws.Range("K$row").Value2 = $homeOpenDt
    $ws.Range("L$row").Value2 = $homeCloseOdds
    $ws.Range("M$row").Value2 = $homeCloseDt

    $ws.Range("N$row").Value2 = $drawOpenOdds
    $ws.Range("O$row").Value2 = $drawOpenDt
    $ws.Range("P$row").Value2 = $drawCloseOdds
    $ws.Range("Q$row").Value2 = $drawCloseDt

    $ws.Range("R$row").Value2 = $awayOpenOdds
    $ws.Range("S$row").Value2 = $awayOpenDt
    $ws.Range("T$row").Value2 = $awayCloseOdds
    $ws.Range("U$row").Value2 = $awayCloseDt

    $ws.Range("V$row").Value2 = $url
}

Add-MatchRow 131 130 45223.91666666666 `
    "Union de Santa Fe" 0 "Defensa y Justicia" 0 `
    2 "21/10/2023 01:42" 2.38 "24/10/2023 21:57" `
    3.36 "21/10/2023 01:42" 3.19 "24/10/2023 21:57" `
    3.78 "21/10/2023 01:42" 3.34 "24/10/2023 21:57" `
    "https://www.betexplorer.com/football/argentina/copa-de-la-liga-profesional/union-de-santa-fe-defensa-y-justicia/4b0f28UG/"

Add-MatchRow 132 131 45224 `
    "Racing Club" 2 "Boca Juniors" 1 `
    2.35 "21/10/2023 01:42" 1.76 "24/10/2023 23:59" `
    3 "21/10/2023 01:42" 3.77 "24/10/2023 23:59" `
    3.51 "21/10/2023 01:42" 4.94 "24/10/2023 23:59" `
    "https://www.betexplorer.com/football/argentina/copa-de-la-liga-profesional/racing-club-boca-juniors/nNka1SqN/"

Add-MatchRow 133 132 45224.10416666666 `
    "Instituto" 0 "Rosario Central" 0 `
    2.08 "20/10/2023 23:12" 1.95 "25/10/2023 02:29" `
    3.15 "20/10/2023 23:12" 3.15 "25/10/2023 02:26" `
    4.07 "20/10/2023 23:12" 4.84 "25/10/2023 02:29" `
    "https://www.betexplorer.com/football/argentina/copa-de-la-liga-profesional/instituto-rosario-central/Isbr5Apb/"

Add-MatchRow 134 133 45224.10416666666 `
    "Argentinos Jrs" 1 "Huracan" 2 `
    1.88 "20/10/2023 23:12" 1.9 "25/10/2023 02:23" `
    3.23 "20/10/2023 23:12" 3.11 "25/10/2023 02:23" `
    4.51 "20/10/2023 23:12" 5.33 "25/10/2023 02:27" `
    "https://www.betexplorer.com/football/argentina/copa-de-la-liga-profesional/argentinos-jrs-huracan/zJ8w6jVi/"
